# Update "想去人数" (want-to-go count) figures in the F column, refreshed
# data pulled at commit a3196b5. The same rows are mirrored on both the
# "展览" sheet and the "全部类型" sheet (which duplicates 展览's rows), so
# every cell is updated on both.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 7494
    "F6"  = 25
    "F7"  = 16
    "F9"  = 5548
    "F12" = 13
    "F13" = 1727
    "F15" = 1116
    "F17" = 5488
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
